# Update the cryptocurrency price and volume data on the active worksheet
# to reflect the refreshed values from the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.180.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +6.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.647.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +16.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.94"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.37"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.642.77"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.64"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.93"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +12.27%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.258.78"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +16.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "71.125.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.644.64"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +16.55%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "516.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +17.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.744"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.59"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.22%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.70%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.55"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.15"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.75"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +12.53%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.78"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.43%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.13"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.10%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.01"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.47%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.347"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.03%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.14"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.16"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.133.85"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "417.32"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0369"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.01%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +12.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +10.58%  "
